$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns X (info_employment_1996) and Y (info_employment_2001) ---
$ws.Range("X1").Value = "info_employment_1996"
$ws.Range("Y1").Value = "info_employment_2001"

$data = @{
    2  = @(33237, 36125)
    3  = @(4672, 6629)
    4  = @(40235, 59085)
    5  = @(20389, 24360)
    6  = @(405370, 556074)
    7  = @(69007, 119642)
    8  = @(43923, 50728)
    9  = @(5561, 8510)
    10 = @(20705, 23298)
    11 = @(141490, 187020)
    12 = @(108047, 146271)
    13 = @(9378, 9027)
    14 = @(8777, 12452)
    15 = @(125353, 147684)
    16 = @(47280, 49394)
    17 = @(31656, 43768)
    18 = @(30855, 46857)
    19 = @(26385, 31116)
    20 = @(25293, 31181)
    21 = @(9413, 12115)
    22 = @(54482, 72232)
    23 = @(103985, 137939)
    24 = @(88150, 101086)
    25 = @(56416, 70510)
    26 = @(13946, 17640)
    27 = @(71709, 86073)
    28 = @(6455, 9845)
    29 = @(25870, 26417)
    30 = @(9892, 21174)
    31 = @(11925, 15468)
    32 = @(116864, 122299)
    33 = @(11751, 15763)
    34 = @(259596, 307699)
    35 = @(58405, 81817)
    36 = @(6195, 8250)
    37 = @(103230, 115124)
    38 = @(25039, 36870)
    39 = @(27398, 37968)
    40 = @(111075, 142171)
    41 = @(10345, 10947)
    42 = @(23027, 30602)
    43 = @(5932, 8383)
    44 = @(41122, 54014)
    45 = @(199114, 282047)
    46 = @(22803, 31718)
    47 = @(6300, 7657)
    48 = @(86044, 145412)
    49 = @(59791, 88203)
    50 = @(12511, 13965)
    51 = @(46060, 54250)
    52 = @(3444, 3936)
}

foreach ($r in 2..52) {
    $vals = $data[$r]

    $cx = $ws.Range("X$r")
    $cx.Value = $vals[0]
    $cx.Font.Name = "Helvetica Neue"
    $cx.Font.Size = 14
    $cx.Font.Color = 2562065
    $cx.NumberFormat = "#,##0"

    $cy = $ws.Range("Y$r")
    $cy.Value = $vals[1]
    $cy.Font.Name = "Helvetica Neue"
    $cy.Font.Size = 14
    $cy.Font.Color = 2562065
    $cy.NumberFormat = "#,##0"

    $ws.Rows($r).RowHeight = 18
}

# Column widths for X and Y (19.1 "characters" renders as the OOXML width=20
# used by the target file, matching its bestFit/customWidth columns)
$ws.Columns("X:Y").ColumnWidth = 19.1

# --- Clear out the old "U.S." aggregate row (row 53), keep styles ---
$ws.Range("A53:W53").ClearContents()

# --- View state ---
$ws.Range("V1").Select()
